# shiny/shinydata/file_info.xlsx — "add duckdb, split app into ui and server"
#
# Spreadsheet-relevant part of that commit:
#  1. The two "chukchi" rows (2019 Chuckchi_Sea_AK samples) were removed from
#     the acoustic_indices sheet, leaving only the "keywest" rows.
#  2. A stray artifact (`"` + trailing spaces, left over from a bad paste) was
#     cleaned out of the 2020_02 / 1024 file-path cell.
#  3. The fish_data and meta tabs were swapped, so the sheet order became
#     acoustic_indices, fish_data, meta (was acoustic_indices, meta, fish_data).
#  4. acoustic_indices became the active tab/selection (C2); meta and
#     fish_data's selections/active-tab flags moved accordingly.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: clean up acoustic_indices -------------------------------------
$ai = $wb.Worksheets.Item("acoustic_indices")

# Fix the trailing `"        ` artifact on the 2020_02/1024 row (row 4) before
# the chukchi rows above it are removed.
$ai.Range("G4").Value = "shinydata/fromLiz/FWRI_KeyWest/16kHz_Decimated/Indices_Native_30sec/2020_02/Decimated_16kHzSR_1024/Acoustic_Indices_16kHz_30sec_1024_20_02.csv"

# Drop the two chukchi/2019 rows entirely (rows 2:3); the keywest rows shift
# up to become rows 2:5.
$ai.Rows("2:3").Delete()

# --- 3: reorder the tabs: fish_data moves before meta ---------------------
$fish = $wb.Worksheets.Item("fish_data")
$meta = $wb.Worksheets.Item("meta")
$fish.Move($meta)

# Re-resolve the worksheet objects by name now that the tab order changed.
$fish = $wb.Worksheets.Item("fish_data")
$meta = $wb.Worksheets.Item("meta")

# --- 4: selections / active tab --------------------------------------------
$fish.Activate()
$fish.Range("D6").Select()

$meta.Activate()
$meta.Range("C7").Select()

$ai.Activate()
$ai.Range("C2").Select()
